$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: integers 1, 2, 3
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

# Row 3: decimals 0.2, 0.3, 1
$ws.Range("A3").Value = 0.2
$ws.Range("B3").Value = 0.3
$ws.Range("C3").Value = 1

# Apply the same style/number format as the rest of the data columns (style index 2 -> 2 decimal numeric format)
$ws.Range("A2:C3").NumberFormat = "0.00"

# Update selection to D3, matching the target workbook view state
$ws.Range("D3").Select()
